$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.432.90'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.371.24'
$ws.Range('E3').Value = '  +3.06%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.66'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.45'
$ws.Range('E6').Value = '  +4.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.519'
$ws.Range('E7').Value = '  -3.27%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.22'
$ws.Range('E10').Value = '  -0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.32'
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.99'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.740.49'
$ws.Range('E15').Value = '  +3.11%  '
$ws.Range('E16').Value = '  +4.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.376.87'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.377.31'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.01'
$ws.Range('E21').Value = '  +3.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0918'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.28'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.73'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.05'
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('E26').Value = '  -0.38%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.82'
$ws.Range('E28').Value = '  +5.97%  '
$ws.Range('E29').Value = '  -4.33%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.29'
$ws.Range('E30').Value = '  +8.05%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.83'
$ws.Range('E31').Value = '  -4.16%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.59'
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '161.96'
$ws.Range('E33').Value = '  -3.38%  '
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.34'
$ws.Range('E36').Value = '  +3.21%  '
$ws.Range('E37').Value = '  +6.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.76'
$ws.Range('E38').Value = '  +12.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.11'
$ws.Range('E39').Value = '  -0.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0744'
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('E41').Value = '  +5.87%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.55'
$ws.Range('E44').Value = '  +11.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.94'
$ws.Range('E45').Value = '  +3.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.007.47'
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('E48').Value = '  +3.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.60'
$ws.Range('E49').Value = '  +7.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '58.20'
$ws.Range('E50').Value = '  +4.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.96'
$ws.Range('E51').Value = '  +1.25%  '
